# fix cargue excel procesos-productos
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header text change: "Monto" -> "Valor" (cell C1, shared-string backed)
$ws.Range("C1").Value = "Valor"

# Hide gridlines on this sheet's view
$excel.ActiveWindow.DisplayGridlines = $false

# Reset selection back to the top-left cell (closest achievable to "no selection override")
$ws.Range("A1").Select()

# Re-layout the columns: add a width for column A and widen B / C
$ws.Columns.Item(1).ColumnWidth = 21
$ws.Columns.Item(2).ColumnWidth = 37.333333333333336
$ws.Columns.Item(3).ColumnWidth = 24.333333333333332
